$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(140, 44235, "4081100", "3011", "Order 4081100 Swish +46721731568", "", "493.75"),
    @(141, 44235, "4081100", "2611", "Order 4081100 Swish +46721731568", "", "59.25"),
    @(142, 44235, "4081100", "1930", "Order 4081100 Swish +46721731568", "553", ""),
    @(143, 44235, "0081101", "3011", "Order 0081101 Swish +46704184802", "", "433.04"),
    @(144, 44235, "0081101", "2611", "Order 0081101Swish +46704184802", "", "51.96"),
    @(145, 44235, "0081101", "1930", "Order 0081101Swish +46704184802", "485", ""),
    @(146, 44236, "4092012", "3011", "Order 4092012 Swish +46703019983", "", "1116.07"),
    @(147, 44236, "4092012", "2611", "Order 4092012 Swish +46703019983", "", "133.93"),
    @(148, 44236, "4092012", "1930", "Order 4092012 Swish +46703019983", "1250", ""),
    @(149, 44236, "4092240", "3011", "Order 4092240 Swish +46763160083", "", "423.21"),
    @(150, 44236, "4092240", "2611", "Order 4092240 Swish +46763160083", "", "50.79"),
    @(151, 44236, "4092240", "1930", "Order 4092240 Swish +46763160083", "474", ""),
    @(152, 44237, "1102020", "3011", "Order 1102020 Swish +46709703734", "", "493.75"),
    @(153, 44237, "1102020", "2611", "Order 1102020 Swish +46709703734", "", "59.25"),
    @(154, 44237, "1102020", "1930", "Order 1102020 Swish +46709703734", "553", ""),
    @(155, 44238, "", "6400", "FACEBK HMUT22KZ62 K6885", "415", ""),
    @(156, 44238, "", "", "FACEBK HMUT22KZ62 K6885", "0", ""),
    @(157, 44238, "", "1930", "FACEBK HMUT22KZ62 K6885", "", "415"),
    @(158, 44238, "", "4010", "SNABBGROSS SOLNA K0135", "845.01", ""),
    @(159, 44238, "", "2645", "SNABBGROSS SOLNA K0135", "101.4", ""),
    @(160, 44238, "", "1930", "SNABBGROSS SOLNA K0135", "", "946.41"),
    @(161, 44239, "0122148", "3011", "Order 0122148 Swish +46732518928", "", "655.36"),
    @(162, 44239, "0122148", "2611", "Order 0122148 Swish +46732518928", "", "78.64"),
    @(163, 44239, "0122148", "1930", "Order 0122148 Swish +46732518928", "734", ""),
    @(164, 44240, "", "5670", "ST1 V#LLINGBY K6885", "645.16", ""),
    @(165, 44240, "", "2641", "ST1 V#LLINGBY K6885", "161.29", ""),
    @(166, 44240, "", "1930", "ST1 V#LLINGBY K6885", "", "806.45"),
    @(167, 44240, "", "4010", "NGROCERIES AB K0135", "272.32", ""),
    @(168, 44240, "", "2645", "NGROCERIES AB K0135", "32.68", ""),
    @(169, 44240, "", "1930", "NGROCERIES AB K0135", "", "305"),
    @(170, 44241, "", "4010", "WILLYS RISSNE K0135", "285.09", ""),
    @(171, 44241, "", "2645", "WILLYS RISSNE K0135", "34.21", ""),
    @(172, 44241, "", "1930", "WILLYS RISSNE K0135", "", "319.3")
)

foreach ($row in $rows) {
    $r = $row[0]
    $dateVal = $row[1]
    $receipt = $row[2]
    $konto = $row[3]
    $descr = $row[4]
    $debet = $row[5]
    $kredit = $row[6]

    # Column A: date serial, formatted like existing rows (style index 2 / YYYY-MM-DD HH:MM:SS)
    $aCell = $ws.Cells.Item($r, 1)
    $aCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $aCell.Value = [double]$dateVal

    # Column B: Receipt Number - may be blank, numeric, or zero-padded text
    $bCell = $ws.Cells.Item($r, 2)
    if ($receipt -eq "") {
        # leave blank - no-op, matches source (cell absent / empty)
    }
    elseif ($receipt.Length -gt 1 -and $receipt.StartsWith("0")) {
        $bCell.NumberFormat = "@"
        $bCell.Value = $receipt
        $bCell.Style = "Normal"
    }
    else {
        $bCell.Value = [double]$receipt
    }

    # Column C: Konto - numeric, may be blank
    $cCell = $ws.Cells.Item($r, 3)
    if ($konto -ne "") {
        $cCell.Value = [double]$konto
    }

    # Column D: Beskrivning - text
    $ws.Cells.Item($r, 4).Value = $descr

    # Column E: Debet - numeric, may be blank
    if ($debet -ne "") {
        $ws.Cells.Item($r, 5).Value = [double]$debet
    }

    # Column F: Kredit - numeric, may be blank
    if ($kredit -ne "") {
        $ws.Cells.Item($r, 6).Value = [double]$kredit
    }
}
